$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nl = [char]10

# --- Helper: write a numeric-looking value as TEXT without polluting styles ---
# Plain Range.Value = "1234" gets auto-coerced to a number by the engine (like
# real Excel). To store it as text (matching the source data, e.g. NUMAR_BILET
# column already stores numeric-looking ids as text) without adding a
# "quote-prefixed" / custom number-format style to styles.xml, we stage the
# text in a scratch cell far outside the used range via a TEXT() formula
# (guaranteed string result), Copy it (copies the *value*, typed as text) into
# the destination, then clear the scratch cell.
function Set-TextValue($cell, $text) {
    $scratch = $ws.Range("AZ1")
    $scratch.Formula = "=TEXT(""$text"",""@"")"
    $scratch.Copy($ws.Range($cell))
    $scratch.Clear()
}

# --- W1 header: style s="1" -> s="3" (same border/fill/font as the other
# header cells, copy formatting from V1 which already carries style 3) ---
$ws.Range("V1").Copy()
$ws.Range("W1").PasteSpecial(-4122)

# --- New column width for W (column 23) ---
# ColumnWidth (COM, character units) = stored OOXML width - 0.83
$ws.Columns.Item(23).ColumnWidth = 13.17

# --- New cell W2 = "1234" (text) ---
Set-TextValue "W2" "1234"

# --- New row 18 ---
$ws.Range("A17").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = 16

$ws.Range("B18").Value = "16-01-2025"
$ws.Range("C18").Value = "DXFDS"
$ws.Range("D18").Value = "SDFDSFDS"
Set-TextValue "E18" "1900106226805"
$ws.Range("I18").Value = "OT-OLT"
$ws.Range("J18").Value = "YES"
$ws.Range("K18").Value = "Elev/Student"
$ws.Range("L18").Value = "YES"
Set-TextValue "M18" "2368"
$ws.Range("N18").Value = $nl
$ws.Range("O18").Value = "NO"
$ws.Range("P18").Value = "NON-APLICABIL"
$ws.Range("Q18").Value = "NON-APLICABIL"
$ws.Range("R18").Value = "NON-APLICABIL"
$ws.Range("S18").Value = "NON-APLICABIL"
$ws.Range("T18").Value = "NO"
$ws.Range("U18").Value = "NON-APLICABIL" + $nl
$ws.Range("V18").Value = "SDFDSFDSFDSFDS" + $nl
Set-TextValue "W18" "9999"

# Assigning a bare "\n" into N18/U18/V18 makes the engine auto-grow the row
# height (customHeight), unlike the source file where every row is left at
# the sheet's default height. AutoFit puts it back to a plain, non-custom row.
$ws.Rows.Item(18).AutoFit()

# --- Extend autofilter range to the new last row ---
$ws.AutoFilterMode = $false
$null = $ws.Range("A1:W18").AutoFilter()

# --- Update the hidden _FilterDatabase defined name to match ---
$names = $wb.Names
$n = $names.Item(1)
$n.RefersTo = "='REGISTRU'!`$A`$1:`$W`$18"
